# Updated remaining queries for C3DC
# Replaces the old `std.id` / `prt.id` style JOIN conditions with the
# fully-qualified `study_id` / `participant_id` column names across every
# SQL query stored on Sheet1 (TabQuery column B and StatQuery column C),
# then restores the view/selection and resizes column C to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-Query {
    param(
        [string]$CellAddress
    )

    $text = $ws.Range($CellAddress).Value2

    $text = $text.Replace(
        'df_participant prt ON std.id = prt."study.id"',
        'df_participant prt ON std.study_id = prt."study.study_id"'
    )
    $text = $text.Replace(
        'df_diagnoses dgn ON prt.id = dgn."participant.id"',
        'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"'
    )
    $text = $text.Replace(
        'df_treatments trt ON prt.id = trt."participant.id"',
        'df_treatments trt ON prt.participant_id = trt."participant.participant_id"'
    )
    $text = $text.Replace(
        'df_treatment_resp trr ON prt.id = trr."participant.id"',
        'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"'
    )
    $text = $text.Replace(
        'df_survival srv ON prt.id = srv."participant.id"',
        'df_survival srv ON prt.participant_id = srv."participant.participant_id"'
    )
    $text = $text.Replace(
        'df_reference_files rfs ON std.id = rfs."study.id"',
        'df_reference_files rfs ON std.study_id = rfs."study.study_id"'
    )

    $ws.Range($CellAddress).Value = $text
}

# StatQuery (column C) + TabQuery (column B) for every tab row.
Update-Query "C2"
Update-Query "B2"
Update-Query "B3"
Update-Query "B4"
Update-Query "B5"
Update-Query "B6"
Update-Query "B7"

# Widen column C (StatQuery) to fit the longer query text, and drop the
# earlier bestFit sizing in favor of an explicit width.
$ws.Columns.Item(3).ColumnWidth = 67.25

# Reset the view: selection moves to B2 and the saved scroll position
# (top-left cell) is cleared.
$ws.Range("B2").Select() | Out-Null
